# Refresh the crypto price/volume snapshot (scheduled GitHub Actions update).
# Coin rows keep their live Price (col D) / Volume(1h) (col E) text; two rows
# (Stellar/Hedera) also swapped rank order, so their Coin name + Link cells
# are updated too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.496.01'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.634.70'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D4').Value = '''0.9992'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''0.9999'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '''305.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('D7').Value = '''0.3754'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '''0.3677'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = '''51.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').Value = '''0.08216'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('D11').Value = '''1.234'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.60%  '
$ws.Range('D12').Value = '''0.9996'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '''22.67'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').Value = '''6.588'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.26%  '
$ws.Range('D15').Value = '''0.00001255'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = '''7.299'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').Value = '1.635.98'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '''94.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').Value = '''0.06979'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '''17.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('E21').Value = '  -1.62%  '
$ws.Range('D22').Value = '''0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '''12.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').Value = '23.499.32'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '''3.184'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.07%  '
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('D27').Value = '''21.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').Value = '''150.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').Value = '''5.323'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').Value = '''134.87'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').Value = '1.817.32'
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').Value = '''2.274'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.90%  '
$ws.Range('D33').Value = '''6.859'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('D34').Value = '''1.027'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.85%  '
$ws.Range('E35').Value = '  +5.08%  '
$ws.Range('D36').Value = '''0.02798'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').Value = '''0.2540'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '''0.08796'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('D39').Value = '''6.094'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '''0.07157'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('D41').Value = '''0.7091'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').Value = '''1.354'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.39%  '
$ws.Range('D43').Value = '''16.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').Value = '''12.36'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('D45').Value = '''0.6575'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').Value = '''2.343'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').Value = '''0.9994'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = '''4.003'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').Value = '''1.213'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '''125.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.45%  '
